# Root_dmg_larvae.xlsx edit script
# - Rename worksheet
# - Translate header row (and, by cascading, the Excel Table column names)
# - Clear 4 mis-recorded "NA" values in column F (rows 132,133,135,136)
# - Update selection / active cell
# - Resize columns B,C,D,E,G (AutoFit-style, bestFit flag dropped)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab
$ws.Name = "Root_dmg (+larvae)"

# 2. Translate the header row (row 1). Because column A:G is bound to the
#    Excel Table ("Tabla53"), updating the header cells also renames the
#    corresponding table columns.
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("D1").Value = "Repeat"
$ws.Range("E1").Value = "Root_weight"
$ws.Range("F1").Value = "N_larvae"
$ws.Range("G1").Value = "Observations"

# 3. These four rows had their larvae count written against the wrong
#    weight record ("NA"); the corrected sheet simply leaves the cell blank.
$ws.Range("F132").ClearContents()
$ws.Range("F133").ClearContents()
$ws.Range("F135").ClearContents()
$ws.Range("F136").ClearContents()

# 4. Resize the columns that no longer need to fit the (now shorter) English
#    labels. (ColumnWidth values chosen so the persisted <col width=.../>
#    lands as close as possible to the authored width; bestFit is dropped
#    automatically once a column is set programmatically.)
$ws.Columns.Item(2).ColumnWidth = 7.799479166666667
$ws.Columns.Item(3).ColumnWidth = 11.619791666666666
$ws.Columns.Item(4).ColumnWidth = 9.072916666666666
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(7).ColumnWidth = 74.70963541666667

# 5. Restore the cursor/selection to the cell the author finished on.
$ws.Range("F133").Select()
